# AirlineRM.xlsx — unit 8 linear optimization lecture file update
#
# This applies the Solver re-run described by the diff: the decision
# variables (Regular seats E5, Discount seats E6) changed from 50/116 to
# 150/100, which ripples through the objective (B8) and the constraint
# block (B11:D15). It also nudges the active-cell selection on Sheet1 and
# (best-effort) stamps VBA CodeNames to mirror the authoring session.

$wb = $excel.ActiveWorkbook

# Best-effort: stamp VBA CodeNames (ThisWorkbook / SheetN) to mirror the
# authoring session that first attached a VBA project to this workbook.
try { $wb.CodeName = "ThisWorkbook" } catch {}
foreach ($sht in $wb.Worksheets) {
    try { $sht.CodeName = $sht.Name } catch {}
}

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Decision variables (Solver "by changing cells") ---
# Regular seats sold
$ws.Range("E5").Value = 150
# Discount seats sold
$ws.Range("E6").Value = 100

# Regular demand upper bound shown alongside the decision (plain data, not
# a formula) also moved to 150 in this run.
$ws.Range("D5").Value = 150

# --- Constraints block ---
# Capacity RHS (plain value, mirrors B11 = E5+E6 which recalculates on its
# own from the changed decision cells above).
$ws.Range("D11").Value = 250

# B8 (objective), B12:B15 and D12:D13 are formulas that recalculate
# automatically from E5/E6/D5/D6, so nothing further to set there.

# --- Selection bookmark left by the editing session ---
$ws.Range("E15").Select()
